# BaoCao_SystemTest_FullCart.xlsx -- "all source and test"
#
# Updates the existing "Quy trinh Mua hang tron ven" row with refreshed
# copy (adds a "kiem tra tong tien" step + clarifies the empty-cart wording),
# then appends two brand-new test-case rows:
#   ST_CART_INV_QTY   - qty-too-low guard rail
#   ST_CART_LARGE_QTY - qty-too-high / stock guard rail
# and widens the columns to fit the new, longer copy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = "`n"

# ---- Row 2: refresh the existing ST_CART_FULL scenario ----------------
$ws.Range("B2").Value = "Quy trình Mua hàng trọn vẹn (CRUD giỏ hàng)"
$ws.Range("C2").Value = "1. Home -> Xem chi tiết SP" + $nl + "2. Chọn Size -> Thêm vào giỏ" + $nl + "3. Vào giỏ -> Update SL lên 2" + $nl + "4. Kiểm tra tổng tiền" + $nl + "5. Xóa SP"
$ws.Range("E2").Value = "Thêm thành công, tính tổng tiền đúng (Price * 2), giỏ hàng trống sau khi xóa"
$ws.Range("F2").Value = "Xong luồng. Trạng thái Giỏ hàng trống: true"

# ---- Row 3: new scenario ST_CART_INV_QTY -------------------------------
$ws.Range("A3").Value = "ST_CART_INV_QTY"
$ws.Range("B3").Value = "Cập nhật số lượng về 0 (Kiểm tra ràng buộc tối thiểu)"
$ws.Range("C3").Value = "1. Đảm bảo giỏ có hàng" + $nl + "2. Nhập số lượng 0" + $nl + "3. Check Alert JS và giá trị input"
$ws.Range("D3").Value = "Quantity: 0"
$ws.Range("E3").Value = "Hiện Browser Alert và số lượng trong input tự động reset về 1"
$ws.Range("F3").Value = "Có Alert: Số lượng tối thiểu là 1! Nếu bạn muốn xóa sản phẩm, vui lòng bấm nút Xóa (thùng rác). | Value reset về: 1"
$ws.Range("G3").Value = "PASS"

# ---- Row 4: new scenario ST_CART_LARGE_QTY -----------------------------
$ws.Range("A4").Value = "ST_CART_LARGE_QTY"
$ws.Range("B4").Value = "Cập nhật số lượng lớn (Kiểm tra ràng buộc tối đa/Tồn kho)"
$ws.Range("C4").Value = "1. Đảm bảo giỏ có hàng" + $nl + "2. Nhập số lượng lớn (9999)" + $nl + "3. Check thông báo HTML và giá trị input"
$ws.Range("D4").Value = "Quantity: 9999"
$ws.Range("E4").Value = "Input reset về Max Stock & Hiện thông báo HTML về tồn kho"
$ws.Range("F4").Value = "Giá trị sau khi nhập: 100 | Alert: Thông báo: Rất tiếc! Sản phẩm này chỉ còn 100 cái trong kho."
$ws.Range("G4").Value = "PASS"

# Give the new PASS cells (G3/G4) the same formatting Excel already uses
# on G2's "PASS" (bold/green) by cloning its format onto the new cells.
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("G4").PasteSpecial(-4122)

# ---- Widen the columns to fit the longer text (bestFit-style autosize) -
# (values chosen so Excel's internal 1/6-character rounding lands as close
# as possible to the POI-computed bestFit widths from the target file)
$ws.Columns.Item(1).ColumnWidth = 19.833333333333332
$ws.Columns.Item(2).ColumnWidth = 52.5
$ws.Columns.Item(3).ColumnWidth = 36.5
$ws.Columns.Item(4).ColumnWidth = 13.5
$ws.Columns.Item(5).ColumnWidth = 67.83333333333333
$ws.Columns.Item(6).ColumnWidth = 103.16666666666667
